$d = $word.ActiveDocument

# Helper: replace the trailing text of a paragraph while preserving any
# leading <w:tab/> run elements untouched (Range.Text / Find-Replace on
# this runtime collapses a leading w:tab into a literal tab character
# inside the new w:t run, which does not match the target XML structure).
function Replace-ParagraphTrailingText {
    param(
        [string]$oldText,
        [string]$newText
    )
    foreach ($para in $d.Paragraphs) {
        $full = $para.Range.Text
        $trimmed = $full.TrimEnd([char]13, [char]7)
        if ($trimmed.EndsWith($oldText)) {
            $tabCount = 0
            while ($tabCount -lt $trimmed.Length -and $trimmed[$tabCount] -eq [char]9) {
                $tabCount++
            }
            $tabsXml = ""
            for ($i = 0; $i -lt $tabCount; $i++) {
                $tabsXml += "<w:tab/>"
            }
            $rng = $para.Range.Duplicate
            $rng.End = $para.Range.End - 1
            $xml = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr>$tabsXml<w:t xml:space=`"preserve`">$newText</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
            $insertStart = $rng.Start
            $rng.Delete()
            $insertPoint = $d.Range($insertStart, $insertStart)
            $insertPoint.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# 1. Remove "(non) " from the precondition text
$d.Content.Find.Execute(
    "Le système a calculé une tournée de livraison (non) valide",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le système a calculé une tournée de livraison valide",
    2) | Out-Null

# 2. Extend the "modifie la plage horaire" bullet
$d.Content.Find.Execute(
    "Le système modifie la plage horaire du point de livraison",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le système modifie la plage horaire du point de livraison et met à jour toutes les heures d’arrivées des points de livraisons suivants",
    2) | Out-Null

# 3. Delete the bullet paragraph "Le système calcule une nouvelle tournée de livraison"
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Le système calcule une nouvelle tournée de livraison") {
        $para.Range.Delete()
        break
    }
}

# 4. Update the confirmation bullet
$d.Content.Find.Execute(
    "L’utilisateur confirme le choix de la nouvelle tournée calculée",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "L’utilisateur confirme le choix de la modification effectuée",
    2) | Out-Null

# 5. Update the "4a." alternate scenario heading -> "3a." (this paragraph
#    starts with a single <w:tab/>, so use the tab-preserving helper)
Replace-ParagraphTrailingText `
    "4a. La nouvelle tournée calculée ne respecte pas les contraintes de l’une tournée de livraison." `
    "3a. La modification ne respecte pas la contrainte des plages horaires" | Out-Null

# 6. Update the corresponding system reaction line (starts with two tabs)
Replace-ParagraphTrailingText `
    "Le système indique que la nouvelle tournée calculée ne respecte les contraintes horaires et retourne à l’étape 1" `
    "Le système met en surbrillance les plages horaires non valide" | Out-Null

# 7. Renumber "5a" to "4a" (starts with a single tab)
Replace-ParagraphTrailingText `
    "5a: L’utilisateur décide de ne pas confirmer son choix" `
    "4a: L’utilisateur décide de ne pas confirmer son choix" | Out-Null

# 8. Renumber "1-5a" to "1-4a"
$d.Content.Find.Execute(
    "1-5a: L’utilisateur indique au système qu’il souhaite annuler la modification de la plage horaire d’un point de livraison",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1-4a: L’utilisateur indique au système qu’il souhaite annuler la modification de la plage horaire d’un point de livraison",
    2) | Out-Null
